$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Tareas": mark several pending tasks as completed ("x") and append
# new tasks (rows 49-52), renaming two of the older abbreviated task notes
# into their full question text (which also frees up their old shared
# string slots).
# ---------------------------------------------------------------------------
$tareas = $wb.Worksheets.Item("Tareas")

# Mark existing tasks as completed
$tareas.Range("D44").Value = "x"
$tareas.Range("D46").Value = "x"
$tareas.Range("D47").Value = "x"
$tareas.Range("D48").Value = "x"

# New task in row 49
$tareas.Range("B49").Value = "Como funciona un flip flop rs con compuertas nand"
$tareas.Range("C49").Value = "10/30/2022"
$tareas.Range("C49").NumberFormat = "dd\-mm\-yy;@"
$tareas.Range("D49").Value = "x"

# New task in row 50
$tareas.Range("B50").Value = "Como se implementa un tren de pulsos con un 555"
$tareas.Range("C50").Value = "10/30/2022"
$tareas.Range("C50").NumberFormat = "dd\-mm\-yy;@"
$tareas.Range("D50").Value = "x"

# Row 51 - reword old abbreviated note into the full task description,
# update its date, mark completed, and underline the index cell (A51) to
# mark the section boundary.
$tareas.Range("B51").Value = "Oscilador con compertas not"
$tareas.Range("C51").Value = "10/30/2022"
$tareas.Range("D51").Value = "x"
$tareas.Range("A51").Font.Underline = $true

# Row 52 - new index number plus reworded note/date/completion mark.
$tareas.Range("A52").Value = 51
$tareas.Range("B52").Value = "¿Qué significa la JK en el flip flop JK?"
$tareas.Range("C52").Value = "10/30/2022"
$tareas.Range("D52").Value = "x"

# ---------------------------------------------------------------------------
# Sheet "Proyectos": add a new project entry in row 6.
# ---------------------------------------------------------------------------
$proyectos = $wb.Worksheets.Item("Proyectos")
$proyectos.Range("B6").Value = "10/30/2022"
$proyectos.Range("B6").NumberFormat = "d-mmm"
$proyectos.Range("C6").Value = "Diseñar un multiplicador x5 de un dígito en BCD.Se desea implementación mínima"

# ---------------------------------------------------------------------------
# Restore on-screen selections for both touched sheets. "Proyectos" is
# selected first so that the final active sheet/tab stays "Tareas", matching
# the saved workbook state.
# ---------------------------------------------------------------------------
$proyectos.Range("C11").Select() | Out-Null
$tareas.Range("D53").Select() | Out-Null
